# Apply the update described in the commit: rename the "Requested quantity"
# headers on the existing sheets and add a new "PO Forecast" sheet with
# forecasted PO quantities.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename header cells in the existing two sheets.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Replicate the bold/bordered header style and the date number-format style
# used on the other sheets so the new sheet reuses the same cell styles.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A22").PasteSpecial(-4122)

# Header row.
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 1; $c -le 4; $c++) {
    $wsForecast.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows.
$aVals = @(45410.99999999999, 45417.99999999999, 45452.99999999999, 45522.99999999999, 45550.99999999999, 45557.99999999999, 45564.99999999999, 45571.99999999999, 45578.99999999999, 45585.99999999999, 45592.99999999999, 45599.99999999999, 45606.99999999999, 45613.99999999999, 45620.99999999999, 45627.99999999999, 45634.99999999999, 45641.99999999999, 45648.99999999999, 45655.99999999999, 45662.99999999999)
$bVals = @(11, 10, 10, 8, 7, 7, 7, 6, 6, 6, 6, 6, 5, 5, 5, 5, 5, 4, 4, 4, 4)
$cVals = @(4.480209072210315, 4.426348164135643, 3.495537288913398, 1.934846236243997, 0.791588749014174, 0.6909233062268656, 0.4619263269180682, 0.6122890686421584, 0.4651286503753816, 0.4139432516641209, 0.08521102905916175, -0.2149417574821492, -0.5112372991017289, -0.379653935563614, -1.103184535376143, -1.213721612773368, -1.424683790030275, -1.308395639480938, -1.712959685218292, -1.669836558144749, -2.13284367193281)
$dVals = @(16.56399961904803, 15.96624333259564, 15.39823542192577, 13.28542178411073, 13.05555894959298, 12.84266932315588, 12.44989715191124, 12.47488348322643, 11.78525695830169, 12.2513043546253, 12.31344426656407, 11.46369395563225, 11.6916204602418, 11.32356042275926, 10.76693528041127, 10.36224361507177, 10.45273228632235, 10.08655519742219, 10.06056181761268, 9.885525531376365, 9.907981389108306)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $aVals[$i]
    $wsForecast.Cells.Item($row, 2).Value = $bVals[$i]
    $wsForecast.Cells.Item($row, 3).Value = $cVals[$i]
    $wsForecast.Cells.Item($row, 4).Value = $dVals[$i]
}

$wsForecast.Range("A1").Select()
